$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算 (A1:F21 -> A1:F22)
# Insert a fresh row 2 for 2021/11/17, push the existing history down.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Cells.Item(2,1).Value = "日期：2021/11/17"
$ws1.Cells.Item(2,2).NumberFormat = "@"
$ws1.Cells.Item(2,2).Value = "202112"
$ws1.Cells.Item(2,3).Value = 17764
$ws1.Cells.Item(2,4).Value = 67327
$ws1.Cells.Item(2,5).Value = 217999808
$ws1.Cells.Item(2,6).Value = 17556
# Keep the sheet's trailing blank row alive as it shifts from row 21 to row 22.
$ws1.Rows.Item(22).OutlineLevel = 0

# ---------------------------------------------------------------------------
# Sheet 2: 散戶多空力道 (A1:B16 -> A1:B17)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Cells.Item(2,1).Value = "日期：2021/11/17"
$ws2.Cells.Item(2,2).Value = -0.13

# ---------------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額 (A1:C16 -> A1:C17)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Cells.Item(2,1).Value = "110年11月17日"
$ws3.Cells.Item(2,2).Value = -2.27
$ws3.Cells.Item(2,3).Value = 73.18000000000001

# ---------------------------------------------------------------------------
# Sheet 4: 大盤多空點位 (A1:B15 -> A1:B16)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Cells.Item(2,1).Value = "110年11月17日"
$ws4.Cells.Item(2,2).Value = 17728.81

# ---------------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位 (A1:N15 -> A1:N16)
# Column A holds bare "yyyy/mm/dd" text that Excel would otherwise parse as
# a real date, so force text first to keep it a literal string.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
$ws5.Cells.Item(2,1).NumberFormat = "@"
$ws5.Cells.Item(2,1).Value = "2021/11/17"
$ws5.Cells.Item(2,2).Value = 44634
$ws5.Cells.Item(2,3).Value = 46292
$ws5.Cells.Item(2,4).Value = -4902
$ws5.Cells.Item(2,5).Value = -3921
$ws5.Cells.Item(2,6).Value = 22930
$ws5.Cells.Item(2,7).Value = 38480
$ws5.Cells.Item(2,8).Value = -6746
$ws5.Cells.Item(2,9).Value = -6259
$ws5.Cells.Item(2,10).Value = -15550
$ws5.Cells.Item(2,11).Value = -487
$ws5.Cells.Item(2,12).Value = 1844
$ws5.Cells.Item(2,13).Value = 2338
$ws5.Cells.Item(2,14).Value = -494
